$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49, pushing the existing rows 49:66 down to 50:67
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record
$ws.Cells.Item(49, 1).Value = 4
$ws.Cells.Item(49, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(49, 3).Value = "Los Lagos"
$ws.Cells.Item(49, 4).Value = 44985
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 6).Value = 100112030
$ws.Cells.Item(49, 7).Value = "Poroto granado"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 80
$ws.Cells.Item(49, 11).Value = 37000
$ws.Cells.Item(49, 12).Value = 37000
$ws.Cells.Item(49, 13).Value = 37000
$ws.Cells.Item(49, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 1480
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"
